# Auto-generated edit script: updates market-price-derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across all 8 sheets
# to match refreshed scheduled-runner data.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 21742134
$ws.Range("I41").Value = 892.4286
$ws.Range("J41").Value = 31253926
$ws.Range("K41").Value = 892.4286
$ws.Range("L41").Value = 31253926
$ws.Range("M41").Value = -452.4286
$ws.Range("N41").Value = -31254806
$ws.Range("H51").Value = 6568.7144
$ws.Range("I51").Value = 5870.75
$ws.Range("K51").Value = 5870.75
$ws.Range("M51").Value = -5386.75
$ws.Range("H113").Value = 13559.471
$ws.Range("I113").Value = 13785.429
$ws.Range("K113").Value = 13785.429
$ws.Range("M113").Value = -10531.429
$ws.Range("H132").Value = 6949.0435
$ws.Range("I132").Value = 8275.4
$ws.Range("K132").Value = 24826.2
$ws.Range("M132").Value = -22296.2
$ws.Range("H137").Value = 4709.3125
$ws.Range("I137").Value = 4540.8
$ws.Range("J137").Value = 4785.909
$ws.Range("K137").Value = 13622.4
$ws.Range("L137").Value = 14357.727
$ws.Range("M137").Value = -11072.4
$ws.Range("N137").Value = -19457.727
$ws.Range("H138").Value = 4255.7544
$ws.Range("J138").Value = 4649.2607
$ws.Range("L138").Value = 13947.7821
$ws.Range("N138").Value = -24227.7821

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 75430.92999999999
$ws.Range("I2").Value = 94171.37
$ws.Range("K2").Value = 94171.37
$ws.Range("M2").Value = -94058.37
$ws.Range("H32").Value = 2824.5
$ws.Range("I32").Value = 2601.5
$ws.Range("K32").Value = 2601.5
$ws.Range("M32").Value = -2314.5
$ws.Range("H74").Value = 13352.72
$ws.Range("J74").Value = 3693.25
$ws.Range("L74").Value = 3693.25
$ws.Range("N74").Value = -5441.25
$ws.Range("H77").Value = 13352.72
$ws.Range("J77").Value = 3693.25
$ws.Range("L77").Value = 18466.25
$ws.Range("N77").Value = -27202.25
$ws.Range("H116").Value = 75430.92999999999
$ws.Range("I116").Value = 94171.37
$ws.Range("K116").Value = 94171.37
$ws.Range("M116").Value = -91877.37
$ws.Range("H122").Value = 7626.091
$ws.Range("I122").Value = 7316.1665
$ws.Range("K122").Value = 21948.4995
$ws.Range("M122").Value = -19498.4995
$ws.Range("H133").Value = 2052000
$ws.Range("J133").Value = 2052000
$ws.Range("L133").Value = 2052000
$ws.Range("N133").Value = -2057060

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 75430.92999999999
$ws.Range("I3").Value = 94171.37
$ws.Range("K3").Value = 94171.37
$ws.Range("M3").Value = -94057.37
$ws.Range("H20").Value = 2121.4285
$ws.Range("I20").Value = 1709.6666
$ws.Range("J20").Value = 3150.8333
$ws.Range("K20").Value = 1709.6666
$ws.Range("L20").Value = 3150.8333
$ws.Range("M20").Value = -1462.6666
$ws.Range("N20").Value = -3644.8333
$ws.Range("H70").Value = 186000
$ws.Range("J70").Value = 186000
$ws.Range("L70").Value = 186000
$ws.Range("N70").Value = -186586
$ws.Range("H73").Value = 186000
$ws.Range("J73").Value = 186000
$ws.Range("L73").Value = 186000
$ws.Range("N73").Value = -188028
$ws.Range("H134").Value = 2686.3264
$ws.Range("I134").Value = 1851
$ws.Range("K134").Value = 5553
$ws.Range("M134").Value = -3018

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2990.389
$ws.Range("J31").Value = 3130.9565
$ws.Range("L31").Value = 3130.9565
$ws.Range("N31").Value = -3720.9565
$ws.Range("H34").Value = 2990.389
$ws.Range("J34").Value = 3130.9565
$ws.Range("L34").Value = 3130.9565
$ws.Range("N34").Value = -3534.9565
$ws.Range("H107").Value = 2532.75
$ws.Range("I107").Value = 1626
$ws.Range("J107").Value = 3180.4285
$ws.Range("K107").Value = 1626
$ws.Range("L107").Value = 3180.4285
$ws.Range("M107").Value = 294
$ws.Range("N107").Value = -7020.4285
$ws.Range("H122").Value = 4573.2856
$ws.Range("I122").Value = 2999.5
$ws.Range("K122").Value = 8998.5
$ws.Range("M122").Value = -6548.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 7028
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 7028
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 21084
$ws.Range("M100").Value = ""
$ws.Range("N100").Value = -22706

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 50000
$ws.Range("J42").Value = 50000
$ws.Range("L42").Value = 50000
$ws.Range("N42").Value = -50970
$ws.Range("H97").Value = 3764.4375
$ws.Range("I97").Value = 4824.5
$ws.Range("K97").Value = 4824.5
$ws.Range("M97").Value = -4328.5
$ws.Range("H102").Value = 7778.15
$ws.Range("I102").Value = 10112.214
$ws.Range("J102").Value = 2332
$ws.Range("K102").Value = 10112.214
$ws.Range("L102").Value = 2332
$ws.Range("M102").Value = -8490.214
$ws.Range("N102").Value = -5576
$ws.Range("H113").Value = 774742.3
$ws.Range("I113").Value = 1004574.9
$ws.Range("J113").Value = 8633.666999999999
$ws.Range("K113").Value = 1004574.9
$ws.Range("L113").Value = 8633.666999999999
$ws.Range("M113").Value = -1002404.9
$ws.Range("N113").Value = -12973.667
$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -52350
$ws.Range("H122").Value = 45634
$ws.Range("I122").Value = 51761.953
$ws.Range("K122").Value = 155285.859
$ws.Range("M122").Value = -152835.859
$ws.Range("H132").Value = 628888.7
$ws.Range("I132").Value = 912675.6
$ws.Range("J132").Value = 4557.4
$ws.Range("K132").Value = 2738026.8
$ws.Range("L132").Value = 13672.2
$ws.Range("M132").Value = -2735496.8
$ws.Range("N132").Value = -18732.2
$ws.Range("H136").Value = 88326.60000000001
$ws.Range("J136").Value = 88326.60000000001
$ws.Range("L136").Value = 264979.8
$ws.Range("N136").Value = -270079.8

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 674258.7
$ws.Range("I7").Value = 1118253.4
$ws.Range("K7").Value = 1118253.4
$ws.Range("M7").Value = -1118141.4
$ws.Range("H40").Value = 1672828.9
$ws.Range("I40").Value = 5003236.5
$ws.Range("J40").Value = 7625
$ws.Range("K40").Value = 5003236.5
$ws.Range("L40").Value = 7625
$ws.Range("M40").Value = -5003100.5
$ws.Range("N40").Value = -7897
$ws.Range("H122").Value = 730214.25
$ws.Range("I122").Value = 594100.6
$ws.Range("K122").Value = 1782301.8
$ws.Range("M122").Value = -1779851.8
$ws.Range("H126").Value = 674258.7
$ws.Range("I126").Value = 1118253.4
$ws.Range("K126").Value = 3354760.2
$ws.Range("M126").Value = -3352290.2
$ws.Range("H132").Value = 3905.6843
$ws.Range("I132").Value = 2435.6428
$ws.Range("K132").Value = 7306.928400000001
$ws.Range("M132").Value = -4776.928400000001

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 64996.668
$ws.Range("J133").Value = 64996.668
$ws.Range("L133").Value = 64996.668
$ws.Range("N133").Value = -75116.66800000001
$ws.Range("H141").Value = 56000
$ws.Range("J141").Value = 56000
$ws.Range("L141").Value = 56000
$ws.Range("N141").Value = -66360
